# [IMP] reference in bank statement import
#
# Insert a new "Reference" column between "Label" and "Partner" in the
# bank statement import sample sheet, and fill in sample reference
# values ("Ref 1" / "Ref 2") for the two example rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Partner") - this shifts the old
# C/D/E (Partner/Amount/Currency) columns one place to the right and the
# new column inherits column C's former width.
$ws.Columns.Item(3).Insert()

# Header row
$ws.Cells.Item(1, 3).Value = "Reference"

# Sample data rows
$ws.Cells.Item(2, 3).Value = "Ref 1"
$ws.Cells.Item(3, 3).Value = "Ref 2"
$ws.Cells.Item(3, 3).Style = "Normal"

$ws.Range("C7").Select() | Out-Null
